# Rapport du 23 Septembre 2025
# Apply the edits described by the commit: fix the duplicated G/H promoter-name
# entries for rows 18-22, clear stale "Transport" (L column) carry-over values,
# tweak column G width, update the active selection, and remove the
# now-obsolete trailing block of rows (172-199) together with the shared
# strings that were only referenced there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 18-22: column H duplicated column G's "promoter name" value.
# Column G should hold the real name (previously sitting in H) and the H
# cell should end up empty. (Use .Value2 - .Value reads misbehave on this
# host for round-tripping cell content between ranges.)
$ws.Range("G18").Value2 = $ws.Range("H18").Value2
$ws.Range("G19").Value2 = $ws.Range("H19").Value2
$ws.Range("G20").Value2 = $ws.Range("H20").Value2
$ws.Range("G21").Value2 = $ws.Range("H21").Value2
$ws.Range("G22").Value2 = $ws.Range("H22").Value2
$ws.Range("H18:H22").Value2 = $null

# --- Clear the stale "Transport" values in column L (rows 13-66 and 69-107).
$ws.Range("L13:L66").Value2 = $null
$ws.Range("L69:L107").Value2 = $null

# --- Column G width tweak (20.9296875 -> ~20.3984375 chars).
$ws.Columns("G").ColumnWidth = 19.5

# --- Remove the trailing duplicated block of rows (172-199).
$ws.Rows("172:199").Delete()

# --- Update the saved view/selection state.
$ws.Range("G18:G22").Select()
